$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape.
# Cells whose new text is a plain decimal number are temporarily
# forced to Text format so COM does not coerce them into a Double
# (e.g. "517.10" -> 517.1 / "0.500" -> 0.5), then the style is reset
# back to 'Normal' so no numFmt/style residue is left on the cell.

$ws.Range('D2').Value = '57.650.31'
$ws.Range('E2').Value = '  -4.29%  '
$ws.Range('D3').Value = '3.090.86'
$ws.Range('E3').Value = '  -6.21%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '517.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -7.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Value = '3.090.76'
$ws.Range('E8').Value = '  -6.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.437'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -9.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.105'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -11.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.366'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -9.97%  '
$ws.Range('D13').Value = '3.629.68'
$ws.Range('E13').Value = '  -6.01%  '
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.59'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.18%  '
$ws.Range('D16').Value = '57.753.94'
$ws.Range('E16').Value = '  -4.12%  '
$ws.Range('D17').Value = '3.099.15'
$ws.Range('E17').Value = '  -5.97%  '
$ws.Range('E18').Value = '  -9.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.72'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.70'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -9.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '335.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -10.22%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.500'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -8.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.166'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.91%  '
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').Value = '0.0₃0897'
$ws.Range('E28').Value = '  -12.11%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.35%  '
$ws.Range('E31').Value = '  -1.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.82'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -9.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.70'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.95'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '157.54'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.69'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.02'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -9.11%  '
$ws.Range('E38').Value = '  -12.03%  '
$ws.Range('D39').Value = '3.128.71'
$ws.Range('E39').Value = '  -6.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '40.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0667'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.79'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -10.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.679'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.93%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.84'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.30%  '
$ws.Range('E46').Value = '  -6.42%  '
$ws.Range('D47').Value = '2.244.58'
$ws.Range('E47').Value = '  -3.28%  '
$ws.Range('E48').Value = '  -10.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.04'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.81%  '
$ws.Range('E51').Value = '  -8.26%  '
